# Update "Chiffres COVID-19 Valais" worksheet with the latest daily figures.
# Only the literal input cells (C = new positive cases, E/F/G = hospital
# breakdown figures) need to be written; columns B/H/J/K are volatile
# formulas (driven by TODAY()) that recompute automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Small corrections to previously-entered daily case counts ---
$ws.Range("C254").Value = 546
$ws.Range("C280").Value = 127
$ws.Range("C611").Value = 61
$ws.Range("C624").Value = 90
$ws.Range("C625").Value = 96

# --- New daily rows (626-629) that were previously blank ---

# Row 626 - 2021-11-12
$ws.Range("C626").Value = 130
$ws.Range("E626").Value = 5
$ws.Range("F626").Value = 4
$ws.Range("G626").Value = 15

# Row 627 - 2021-11-13
$ws.Range("C627").Value = 75
$ws.Range("E627").Value = 5
$ws.Range("F627").Value = 4
$ws.Range("G627").Value = 15

# Row 628 - 2021-11-14
$ws.Range("C628").Value = 29
$ws.Range("E628").Value = 5
$ws.Range("F628").Value = 4
$ws.Range("G628").Value = 16

# Row 629 - 2021-11-15
$ws.Range("C629").Value = 8
$ws.Range("E629").Value = 5
$ws.Range("F629").Value = 4
$ws.Range("G629").Value = 16
